$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (H1) onto the new
# header cells so I1/J1 match the existing header style.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Headers
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF)
$data = @{
    2  = @(7, 8)
    3  = @(9, 9)
    4  = @(9, 9)
    5  = @(8, 9)
    6  = @(9, 9)
    7  = @(7, 7)
    8  = @(8, 8)
    9  = @(8, 8)
    10 = @(8, 8)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
